# "Generate Report for Archive"
#
# 1. The localization status moved on from "Ready for handoff" to
#    "In Translation" - update every cell that was showing the old
#    status (Overview!E2:F3, and the Status column - C2:C3 - on each
#    per-language sheet) so they all read the new value.
# 2. Now that the status text is shorter, those status columns are
#    narrowed to fit it again.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Re-fit the status columns to the new, shorter text (~13.41 chars wide).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
